# Update the dSF column (F) values for specific rows to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 0
    4  = 8
    5  = -11
    10 = -7
    12 = 11
    14 = -5
    15 = -9
    16 = 6
    18 = -2
    22 = 1
    29 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
